# Update countries & provincias Spain
#
# The underlying data source refreshed its COVID-19 snapshot. For several
# countries the new "Casos totales" (column B) value changed enough to
# change their sort rank relative to a neighbouring row, so the row's
# entire B:H record (and, where the neighbour swap happened, the country
# name in column A) now belongs to a different row than before:
#   - row 15/16: Canada now outranks Brasil
#   - rows 21-25: Irlanda now outranks India; Peru now outranks Suecia,
#     which now outranks Israel
#   - rows 138-140: Aruba now outranks Guayana Francesa, which now
#     outranks Gabon
# Rows 18, 76, 106 and 141 keep their country but received refreshed
# numbers only.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-PaisRow($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

Set-PaisRow 15  "Canada"           31642 1536 10325 20007 557  115 1310
Set-PaisRow 16  "Brasil"           30961 278  14026 14979 6634 9   1956
Set-PaisRow 18  "Suiza"            27078 346  16400 9353  386  44  1325
Set-PaisRow 21  "Irlanda"          13980 709  77    13373 156  44  530
Set-PaisRow 22  "India"            13835 405  1777  11606 0    4   452
Set-PaisRow 23  "Peru"             13489 998  6120  7069  169  26  300
Set-PaisRow 24  "Suecia"           13216 676  550   11266 482  67  1400
Set-PaisRow 25  "Israel"           12982 224  3126  9705  168  9   151
Set-PaisRow 76  "Lituania"         1149  21   210   906   14   1   33
Set-PaisRow 106 "Jordania"         407   5    265   135   5    0   7
Set-PaisRow 138 "Aruba"            96    1    43    51    1    0   2
Set-PaisRow 139 "Guayana Francesa" 96    10   61    35    2    0   0
Set-PaisRow 140 "Gabon"            95    15   6     88    0    0   1
Set-PaisRow 141 "Monaco"           94    1    20    71    3    0   3
